$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-12 Thursday" "2026-02-13 Friday"

Replace-Text "304×3=" "572×4="
Replace-Text "434×5=" "537×5="
Replace-Text "343×4=" "416×6="
Replace-Text "358×9=" "726×5="
Replace-Text "849×8=" "262×7="

Replace-Text "591×9=" "867×7="
Replace-Text "899×6=" "378×3="
Replace-Text "207×8=" "447×2="
Replace-Text "333×8=" "281×4="
Replace-Text "744×5=" "574×4="

Replace-Text "862×8=" "858×7="
Replace-Text "558×2=" "808×8="
Replace-Text "492×9=" "624×4="
Replace-Text "243×9=" "637×7="
Replace-Text "940×3=" "945×4="

Replace-Text "440×3=" "766×6="
Replace-Text "516×7=" "367×7="
Replace-Text "236×5=" "706×8="
Replace-Text "426×8=" "776×3="
Replace-Text "841×8=" "550×5="

Replace-Text "745×3=" "358×8="
Replace-Text "330×6=" "702×6="
Replace-Text "791×6=" "182×4="
Replace-Text "650×4=" "754×5="
Replace-Text "502×9=" "644×9="
